$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert new column B for "Week_Start_Date" - shifts old B..I to C..J
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "Week_Start_Date"

$weekStarts = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

# Force text format so the date-like strings are NOT auto-converted to date serials
$ws.Range("B2:B17").NumberFormat = "@"
for ($i = 0; $i -lt $weekStarts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStarts[$i]
}

# Fix Week labels in column A: W01 -> W1 ... W16 -> W16 (strip leading zero)
for ($i = 1; $i -le 16; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "W$i"
}

# Updated MyForecast values (now column D after the insert)
$myForecast = @(49,52,63,53,47,50,63,55,51,55,67,60,43,45,55,45)
for ($i = 0; $i -lt $myForecast.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]
}

# is_holiday_week column (now J) becomes boolean typed
for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $false
}

# Update Summary sheet - keep values as plain text (matches existing column formatting)
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9:B13").NumberFormat = "@"
$summary.Range("B9").Value = "853"
$summary.Range("B10").Value = "432"
$summary.Range("B11").Value = "218"
$summary.Range("B12").Value = "67"
$summary.Range("B13").Value = "2025-03-16"
